$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '95.286.47'
$ws.Range("E2").Value = '  +1.69%  '
$ws.Range("D3").Value = '3.610.67'
$ws.Range("E3").Value = '  +3.64%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '235.36'
$ws.Range("E5").Value = '  -0.04%  '
$ws.Range("D6").Value = '656.50'
$ws.Range("E6").Value = '  +5.05%  '
$ws.Range("E7").Value = '  +1.78%  '
$ws.Range("D8").Value = '0.400'
$ws.Range("E8").Value = '  +1.88%  '
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("D10").Value = '0.996'
$ws.Range("E10").Value = '  -0.18%  '
$ws.Range("D11").Value = '3.607.27'
$ws.Range("E11").Value = '  +3.63%  '
$ws.Range("E12").Value = '  +0.11%  '
$ws.Range("D13").Value = '42.00'
$ws.Range("E13").Value = '  -1.89%  '
$ws.Range("D14").Value = '6.37'
$ws.Range("E14").Value = '  +1.79%  '
$ws.Range("D15").Value = '4.310.44'
$ws.Range("E15").Value = '  +3.91%  '
$ws.Range("D16").Value = '95.271.44'
$ws.Range("E16").Value = '  +1.88%  '
$ws.Range("E17").Value = '  +1.19%  '
$ws.Range("D18").Value = '3.615.56'
$ws.Range("E18").Value = '  +3.52%  '
$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").Value = '7.91'
$ws.Range("E19").Value = '  -5.18%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = '12.89'
$ws.Range("E20").Value = '  +3.85%  '
$ws.Range("D21").Value = '17.91'
$ws.Range("E21").Value = '  -0.28%  '
$ws.Range("D22").Value = '3.51'
$ws.Range("E22").Value = '  +3.63%  '
$ws.Range("D23").Value = '507.49'
$ws.Range("E23").Value = '  -1.96%  '
$ws.Range("D24").Value = '0.476'
$ws.Range("E24").Value = '  -4.25%  '
$ws.Range("D25").Value = '0.0000194'
$ws.Range("E25").Value = '  +6.40%  '
$ws.Range("D26").Value = '6.58'
$ws.Range("E26").Value = '  -2.12%  '
$ws.Range("D27").Value = '95.15'
$ws.Range("E27").Value = '  -0.36%  '
$ws.Range("D28").Value = '3.806.64'
$ws.Range("E28").Value = '  +3.63%  '
$ws.Range("D29").Value = '12.44'
$ws.Range("E29").Value = '  +1.90%  '
$ws.Range("D30").Value = '3.06'
$ws.Range("E30").Value = '  +3.41%  '
$ws.Range("B31").Value = 'Dai'
$ws.Range("C31").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  -0.13%  '
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").Value = '11.19'
$ws.Range("E32").Value = '  -1.68%  '
$ws.Range("E33").Value = '  +0.71%  '
$ws.Range("B34").Value = 'Binance-PegBSC-USD'
$ws.Range("C34").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D34").Value = '0.998'
$ws.Range("E34").Value = '  -0.18%  '
$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D35").Value = '32.64'
$ws.Range("E35").Value = '  +9.56%  '
$ws.Range("D36").Value = '0.176'
$ws.Range("E36").Value = '  -1.08%  '
$ws.Range("D37").Value = '0.559'
$ws.Range("E37").Value = '  +0.37%  '
$ws.Range("D38").Value = '8.08'
$ws.Range("E38").Value = '  +7.39%  '
$ws.Range("D39").Value = '563.28'
$ws.Range("E39").Value = '  -0.97%  '
$ws.Range("E40").Value = '  +2.03%  '
$ws.Range("B41").Value = 'USDe'
$ws.Range("C41").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").Value = '0.150'
$ws.Range("E42").Value = '  +1.58%  '
$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D43").Value = '0.914'
$ws.Range("E43").Value = '  -0.58%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = '36.53'
$ws.Range("E44").Value = '  +44.04%  '
$ws.Range("E45").Value = '  +4.77%  '
$ws.Range("E46").Value = '  +4.13%  '
$ws.Range("D47").Value = '23.58'
$ws.Range("E47").Value = '  -0.65%  '
$ws.Range("E48").Value = '  +4.87%  '
$ws.Range("D49").Value = '0.0411'
$ws.Range("E49").Value = '  -1.33%  '
$ws.Range("D50").Value = '3.58'
$ws.Range("E50").Value = '  +0.79%  '
$ws.Range("D51").Value = '53.20'
$ws.Range("E51").Value = '  -0.48%  '
